$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 'https://lh3.googleusercontent.com/pw/AP1GczP43uebQ4gDcWJLejZ7DTWX7sZ1Kve4EdQrf0dyDoz8BzIA7IsnM5TN7OgnN22XYm2PT-aUoxdHoaRZ2-00U0kqfhlR0StybFhWteZQGGkSdI0wgvSpFW-z2t1JsDtX7kG-cA8EwjoiI-gEDucMDlBo1A=w500-h309-s-no-gm?authuser=0'
$ws.Range("C13").Value = 'Gimnasia Pasiva'
$ws.Range("D13").Value = '45.000'
$ws.Range("E13").Value = 'Gimnasia Pasiva 8 Electrodos De Pila Y Conectar'
$ws.Range("F13").Value = 'edward quevedo'
$ws.Range("G13").Value = '3224336216'
$ws.Range("H13").Value = 'no'
$ws.Range("I13").Value = 'https://lh3.googleusercontent.com/pw/AP1GczP43uebQ4gDcWJLejZ7DTWX7sZ1Kve4EdQrf0dyDoz8BzIA7IsnM5TN7OgnN22XYm2PT-aUoxdHoaRZ2-00U0kqfhlR0StybFhWteZQGGkSdI0wgvSpFW-z2t1JsDtX7kG-cA8EwjoiI-gEDucMDlBo1A=w500-h309-s-no-gm?authuser=0'
$ws.Range("J13").Value = 'https://lh3.googleusercontent.com/pw/AP1GczNyOyhYX-cmhTtPph9TaYAmzjNUy75kzfsl_jWyZkOGS0C6TTJBQkx06Yv_3EOMmDj4ULQMr4bSM7_pqkXA0QINGouykld9Mx_Zsl7JOg1ml9g98oquSmlvGknceIZkh7R0TOoBa3BE0CdMZ2unjkNBbg=w432-h500-s-no-gm?authuser=0'
$ws.Range("K13").Value = 'https://lh3.googleusercontent.com/pw/AP1GczNFTrSqJt6NsTEkkAtGM5vvJBbjHyai9N2pi42AVY88GnQ-gHyQKoT-iHh8y2TaQ8qlEaRpzYqYKKYaJIZBfb1kK18W2URs7xJ5_NrKYx7bKhZsN1ufh6oSjdFIKGpH2iFf5M6gYlqPYd0Tp3H-PdqF3g=w500-h489-s-no-gm?authuser=0'
$ws.Range("L13").Value = 'https://lh3.googleusercontent.com/pw/AP1GczNDxVXm5qR802junRd2lyek41EuaONp19hqhF-CFq63rMIKmLkP025CNuWiJNgG6_2NTV8PuN-CO2df_Blvj2bNjNlnQn1fJLX-m4xEtq_Bc3vPmu0vtRUSawxRt-hxgV6lr5UFs10Uv9C8Ixt8ma4ibA=w500-h363-s-no-gm?authuser=0'

$ws.Range("B14").Value = 'https://lh3.googleusercontent.com/pw/AP1GczMThvkBf7W_2DmcdJ0jzASx1GAhBztiFTBvJLHfKAoy2BpG_OoKs1l6RutqtzUR1AslF75pguaxp2eH7tFYlPZflsQXk88FU6o7dUCLuO4Ke1BzOW7t-NVe9t5I7CQj4tYv_CBuiu5ni_24BxQTOxOAug=w500-h428-s-no-gm?authuser=0'
$ws.Range("C14").Value = 'Cepillo Secador 3 En 1'
$ws.Range("D14").Value = '50.000'
$ws.Range("E14").Value = 'Cepillo Secador 3 En 1, Cepilla, Seca, Alisa'
$ws.Range("F14").Value = 'edward quevedo'
$ws.Range("G14").Value = '3224336217'
$ws.Range("H14").Value = 'no'
$ws.Range("I14").Value = 'https://lh3.googleusercontent.com/pw/AP1GczMThvkBf7W_2DmcdJ0jzASx1GAhBztiFTBvJLHfKAoy2BpG_OoKs1l6RutqtzUR1AslF75pguaxp2eH7tFYlPZflsQXk88FU6o7dUCLuO4Ke1BzOW7t-NVe9t5I7CQj4tYv_CBuiu5ni_24BxQTOxOAug=w500-h428-s-no-gm?authuser=0'

$ws.Range("B15").Value = 'https://lh3.googleusercontent.com/pw/AP1GczMFGpnhAYcuxekUHmsBMgL0NoSjeqF2TIBV6-rOwBHYmvl6pKbsd89tg8Nm6_TTvwBVSSETC1Kb5DODQHHIQnyjXt2N06ar8XQng1lftI0gy_97zK34Koem9FuO6UfAuX59BdtGCPftVmY7d0q862ibIg=w794-h633-s-no-gm?authuser=0'
$ws.Range("C15").Value = 'Combo Gamer'
$ws.Range("D15").Value = '60.000'
$ws.Range("E15").Value = 'Combo Gamer Negro RGB Teclado Y Mouse'
$ws.Range("F15").Value = 'edward quevedo'
$ws.Range("G15").Value = '3224336217'
$ws.Range("H15").Value = 'no'
$ws.Range("I15").Value = 'https://lh3.googleusercontent.com/pw/AP1GczMFGpnhAYcuxekUHmsBMgL0NoSjeqF2TIBV6-rOwBHYmvl6pKbsd89tg8Nm6_TTvwBVSSETC1Kb5DODQHHIQnyjXt2N06ar8XQng1lftI0gy_97zK34Koem9FuO6UfAuX59BdtGCPftVmY7d0q862ibIg=w794-h633-s-no-gm?authuser=0'
$ws.Range("J15").Value = 'https://lh3.googleusercontent.com/pw/AP1GczN_YeDAx1aukAWLluS_995Z9heN-qKegOKTLNuPHJMEBQG-qqL7P4AaQtYt8OwUOXlPuRX1PvM1hiZnGCK3hJQkbedX-ycwIxj8mx5KRjVvYlwTFQ8RnVYrNChSOLxJ8bue909bF0WKzI2Skwf3LsstnA=w600-h600-s-no-gm?authuser=0'

$ws.Range("J17").Select()
